$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings, written in the exact order needed so the
# resulting shared-strings table indices line up with the target file.
$filesTab = @'
FilesTab
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
   WHERE c.disease = "Adenocarcinoma of the cervix"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE c.disease = "Adenocarcinoma of the cervix"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

$caseQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
 WHERE c.disease = "Adenocarcinoma of the cervix"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# Row 2 (CasesTab): query in B2 changes to the case-list query;
# the stat query in C2 is replaced with the new 'Trials/Cases/Files' query.
$ws.Range("A3").Value = $filesTab
$ws.Range("C2").Value = $statQuery
$ws.Range("B3").Value = $filesQuery
$ws.Range("B2").Value = $caseQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

# Row heights: row2 grows slightly, new row3 needs the max row height.
$ws.Rows(2).RowHeight = 195
$ws.Rows(3).RowHeight = 409.5

# Wrap text formatting on the query columns (style already used by B2/C2).
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Restore selection to B2, matching the saved workbook view state.
$ws.Range("B2").Select()
